# Fix list level numbering: top-level list items were being emitted one
# level too deep (and nested list items two levels too deep). Decrement
# the indent level of the affected paragraphs by one.

$p = $ppt.ActivePresentation

# Slide 1: "Content Placeholder 2" shape
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(2).TextFrame.TextRange

# "Bullet item with inline code" : lvl 1 -> lvl 0  (IndentLevel 2 -> 1)
$tr1.Paragraphs(2, 1).IndentLevel = 1

# "with nested" : lvl 2 -> lvl 1  (IndentLevel 3 -> 2)
$tr1.Paragraphs(4, 1).IndentLevel = 2

# Slide 2: "Content Placeholder 2" shape
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange

# "Nested" : lvl 1 -> lvl 0  (IndentLevel 2 -> 1)
$tr2.Paragraphs(2, 1).IndentLevel = 1

# Slide 3: "Content Placeholder 2" shape
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange

# "A total alternative for head" : lvl 1 -> lvl 0  (IndentLevel 2 -> 1)
$tr3.Paragraphs(1, 1).IndentLevel = 1
